$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate rows 839-847 with event data ---
$row = 839
$ws.Cells.Item($row, 1).Value = 45968
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'RAVE INDUSTRY'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Stollen134'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Dortmund'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DPwcvJGDX7Q/?igsh=MnFuYzFtMnFyeDho', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DPwcvJGDX7Q/?igsh=MnFuYzFtMnFyeDho')

$row = 840
$ws.Cells.Item($row, 1).Value = 45946
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'SALON ELECTRONIQUE'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'SNRS'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Dortmund'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/p/DPjOUNdCDnF/?igsh=c3Y1YXMyYmtjMXl6', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/p/DPjOUNdCDnF/?igsh=c3Y1YXMyYmtjMXl6')

$row = 841
$ws.Cells.Item($row, 1).Value = 45955
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'CLUB NIGHT'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'SNRS'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Dortmund'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DPBpSFSjGnL/?igsh=MTR1YndkaDI1NWw5ZQ==', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DPBpSFSjGnL/?igsh=MTR1YndkaDI1NWw5ZQ==')

$row = 842
$ws.Cells.Item($row, 1).Value = 45945
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'ROOT 150 MIN RAVE'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'AREA 15'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Bochum'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DPg3q8xiC7S/?igsh=a240OGU5cWR5ZHYz', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DPg3q8xiC7S/?igsh=a240OGU5cWR5ZHYz')

$row = 843
$ws.Cells.Item($row, 1).Value = 45990
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'POLAAR CLUB FESTIVAL'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Rotunde, Area 15, Klub Kurt'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Bochum'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DPMZsmyCkQI/?igsh=bzU5YWg4M2FiMzQz', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DPMZsmyCkQI/?igsh=bzU5YWg4M2FiMzQz')

$row = 844
$ws.Cells.Item($row, 1).Value = 45961
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'TAKTLOS HALLOWEEN'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Panzergarten und Westend'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Essen'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/p/DPZOzTEDU1c/?igsh=MWJrODkxY2M4aXF4bw==', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/p/DPZOzTEDU1c/?igsh=MWJrODkxY2M4aXF4bw==')

$row = 845
$ws.Cells.Item($row, 1).Value = 45976
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'PUMP'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Fusion'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Münster'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DNOGdh9MSaF/?igsh=MWR3MXdidnJudDViYQ==', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DNOGdh9MSaF/?igsh=MWR3MXdidnJudDViYQ==')

$row = 846
$ws.Cells.Item($row, 1).Value = 45948
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'EHRENKLUB'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Odonien'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Köln'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DPRXPaHiE_9/?igsh=MWVobTlnNmU0NDQ0aA==', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DPRXPaHiE_9/?igsh=MWVobTlnNmU0NDQ0aA==')

$row = 847
$ws.Cells.Item($row, 1).Value = 45976
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = 'EHRENKLUB'
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = 'Schrotty'
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = 'Köln'
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), 'https://www.instagram.com/reel/DOOY05uCEkV/?igsh=MTBqY2F0Z3JkbTRwbg==', [Type]::Missing, [Type]::Missing, 'https://www.instagram.com/reel/DOOY05uCEkV/?igsh=MTBqY2F0Z3JkbTRwbg==')

# --- Extend sheet with 28 new blank rows (852-879), matching existing blank-row style ---
$srcBlank = $ws.Range("A851:E851")
$dstBlank = $ws.Range("A852:E879")
$srcBlank.Copy($dstBlank)
$ws.Range("A852:E879").RowHeight = 15

Write-Output $ws.UsedRange.Address()
